# adding misc Win10 pro defaults
$wb = $excel.ActiveWorkbook

# Rename "User Public root directory" sheet to "User Public root".
# This also updates the corresponding _xlnm._FilterDatabase defined name
# that references the sheet by its (quoted) name.
$wsPublic = $wb.Worksheets.Item("User Public root directory")
$wsPublic.Name = "User Public root"

# Update the selected cell on the (now renamed) "User Public root" sheet
# from A1:E2 to the single cell N27.
$wsPublic.Activate()
$wsPublic.Range("N27").Select()

# Update the selected/active cell on the "User root dir" sheet
# (the tab that remains selected) from D25 to H25.
$wsUserRoot = $wb.Worksheets.Item("User root dir")
$wsUserRoot.Activate()
$wsUserRoot.Range("H25").Select()
